$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9300694554254023
$ws.Range("J2").Value = 0.9300694554254023
$ws.Range("M2").Value = 8.676671000000001
$ws.Range("N2").Value = 26.030013
$ws.Range("O2").Value = 0.1325240072999665
$ws.Range("P2").Value = 0.1325240072999665
$ws.Range("Q2").Value = 21.75780530191467
$ws.Range("R2").Value = 195.820247717232
$ws.Range("S2").Value = 0.1232565313002719
$ws.Range("T2").Value = 0.1232565313002719

# Row 3
$ws.Range("I3").Value = 0.9300694554254023
$ws.Range("J3").Value = 0.9300694554254023
$ws.Range("M3").Value = 37.74750533333334
$ws.Range("O3").Value = 0.5765403197090441
$ws.Range("P3").Value = 0.576540319709044
$ws.Range("Q3").Value = 94.65644965398045
$ws.Range("R3").Value = 851.9080468858241
$ws.Range("S3").Value = 0.5362225411825779
$ws.Range("T3").Value = 0.5362225411825778

# Row 4
$ws.Range("I4").Value = 0.9300694554254023
$ws.Range("J4").Value = 0.9300694554254023
$ws.Range("M4").Value = 19.04827033333333
$ws.Range("N4").Value = 57.144811
$ws.Range("O4").Value = 0.2909356729909895
$ws.Range("P4").Value = 0.2909356729909895
$ws.Range("Q4").Value = 47.76584905096711
$ws.Range("R4").Value = 429.892641458704
$ws.Range("S4").Value = 0.2705903829425525
$ws.Range("T4").Value = 0.2705903829425525

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1885443333333333
$ws.Range("H5").Value = 0.5656329999999999
$ws.Range("I5").Value = 0.06993054457459773
$ws.Range("J5").Value = 0.06993054457459771
$ws.Range("M5").Value = 8.676671000000001
$ws.Range("N5").Value = 26.030013
$ws.Range("O5").Value = 0.1325240072999665
$ws.Range("P5").Value = 0.1325240072999665
$ws.Range("Q5").Value = 1.635937149247667
$ws.Range("R5").Value = 14.723434343229
$ws.Range("S5").Value = 0.009267475999694621
$ws.Range("T5").Value = 0.009267475999694618

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1885443333333333
$ws.Range("H6").Value = 0.5656329999999999
$ws.Range("I6").Value = 0.06993054457459773
$ws.Range("J6").Value = 0.06993054457459771
$ws.Range("M6").Value = 37.74750533333334
$ws.Range("O6").Value = 0.5765403197090441
$ws.Range("P6").Value = 0.576540319709044
$ws.Range("Q6").Value = 7.117078228069778
$ws.Range("R6").Value = 64.05370405262799
$ws.Range("S6").Value = 0.04031777852646613
$ws.Range("T6").Value = 0.04031777852646611

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1885443333333333
$ws.Range("H7").Value = 0.5656329999999999
$ws.Range("I7").Value = 0.06993054457459773
$ws.Range("J7").Value = 0.06993054457459771
$ws.Range("M7").Value = 19.04827033333333
$ws.Range("N7").Value = 57.144811
$ws.Range("O7").Value = 0.2909356729909895
$ws.Range("P7").Value = 0.2909356729909895
$ws.Range("Q7").Value = 3.591443431151444
$ws.Range("R7").Value = 32.322990880363
$ws.Range("S7").Value = 0.02034529004843698
$ws.Range("T7").Value = 0.02034529004843698
